$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: update rpc-reply message-id uuid
$f2 = $ws.Range("F2").Value2
$f2 = $f2 -replace "ab43be2c-7470-41cf-a114-beb491971b0d", "fa681b9d-c96a-421f-a286-3d07d5f51129"
$ws.Range("F2").Value2 = $f2

# G2: update protocol identifier (add xmlns + prefix) and rename BGP_65000 -> default (both occurrences)
$g2 = $ws.Range("G2").Value2
$g2 = $g2 -replace "<identifier>BGP</identifier>", '<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>'
$g2 = $g2 -replace "<name>BGP_65000</name>", "<name>default</name>"
$ws.Range("G2").Value2 = $g2
